# Rename the "SSx" sheets to their new descriptive names
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("SS1").Name = "0,1_0,5_1"
$wb.Worksheets.Item("SS2").Name = "0,1_0,5_0,5"
$wb.Worksheets.Item("SS3").Name = "0,1_0,5_3"
$wb.Worksheets.Item("SS4").Name = "0,5_1_1"
$wb.Worksheets.Item("SS5").Name = "0,05_0,2_1"

# Update the selected ranges / active cells on a few sheets
$wb.Worksheets.Item("0,1_0,5_1").Range("A1:B100").Select()
$wb.Worksheets.Item("0,1_0,5_3").Range("A1:B100").Select()
$wb.Worksheets.Item("0,05_0,2_1").Range("F6").Select()

# Restore the originally active sheet so the workbook-level view is unchanged
$wb.Worksheets.Item(1).Activate()
